$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6 (shifting "Johnson & Johnson..." block and
# everything below it down by one row), then fill it with the new
# alternative vaccine name "Comirnaty" (Pfizer's brand name).
$ws.Rows("6:6").Insert()
$ws.Cells.Item(6, 1).Value = "Comirnaty"

# Match the author's final selection state (cell A7 selected).
$ws.Range("A7").Select()
